$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Journal")
$ws2 = $wb.Worksheets.Item("Totaux")

# 1) Fill in row 10 of the Journal sheet with the new journal entry
#    (new task added after the chef de projet's feedback)
$ws1.Range("A10").Value = 44964
$ws1.Range("B10").Value = 2
$ws1.Range("D10").Value = "Documentation"
$ws1.Range("E10").Value = "Modification de la planification initiale du projet, suite au retour du chef de projet"

# 2) The longer description text no longer fits the old column width -
#    widen column E (mirrors Excel re-running best-fit on the column)
$ws1.Columns.Item(5).ColumnWidth = 75.140625

# 3) Move the selection on the Journal sheet to the next empty row
[void]$ws1.Range("C10").Select()

# 4) Switch the active tab over to "Totaux"
[void]$ws2.Activate()
[void]$ws2.Range("C8").Select()
